$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 575
$ws.Range("I2").Value = 1479
$ws.Range("J2").Value = 6037
$ws.Range("K2").Value = 33
$ws.Range("L2").Value = 1610
$ws.Range("M2").Value = 102
$ws.Range("N2").Value = 1109
$ws.Range("O2").Value = 5
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 78
$ws.Range("S2").Value = 646
$ws.Range("T2").Value = 1065
$ws.Range("U2").Value = 75
$ws.Range("V2").Value = 9148
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 9519
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 145
$ws.Range("AA2").Value = 47
